$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    103.8374204181422,
    117.1151164473825,
    119.4110606755175,
    136.2096301478599,
    158.2918350324003,
    163.1849125657315,
    177.1411248488167,
    193.4043851827645,
    215.4590554986885,
    224.336968030124,
    229.2257191944717,
    243.1499849728491,
    253.5736295707525,
    267.3199243764205,
    270.7222645061136,
    287.3662580521499,
    293.3167470268107,
    306.4229130947654,
    315.9232784034401,
    328.2875258911848,
    344.6166804119313,
    352.7094509487043,
    352.7441020764021,
    361.7056148216471,
    371.587840846849,
    388.6551514251913,
    404.3854989028974,
    407.7086747031861,
    419.0683999097981,
    411.724143966949,
    420.9523327099513,
    423.3908706175966,
    443.4889444088479,
    463.0784582705287,
    463.0165240485661,
    474.4062461289118,
    498.4968046305327,
    496.7905357500559,
    511.9083059793905,
    527.4384052033954,
    539.1244830164649,
    542.022563681335,
    558.2228568778291,
    566.8859578507953,
    573.8101819391982,
    568.5705620517134,
    576.2008181593129,
    596.4006279503204,
    599.6277182022559,
    600.3262026089017,
    609.4329745905666,
    612.0423566192328,
    621.3326023214269,
    639.4546735964685,
    631.2355020184473,
    648.7264335148374,
    649.94556550322,
    650.0099109054011,
    652.4215876579328,
    651.9769151425614,
    654.0717279037372
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
